# Applies the cryptos.xlsx price/volume/coin-order update described in the
# commit's OOXML diff (GitHub Actions crypto-list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "97.224.38"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").Value = "3.726.19"
$ws.Range("E3").Value = "  +1.47%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.68"
$ws.Range("E5").Value = "  -0.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.91"
$ws.Range("E6").Value = "  +0.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "660.11"
$ws.Range("E7").Value = "  +0.51%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.437"
$ws.Range("E8").Value = "  +2.56%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.06"
$ws.Range("E10").Value = "  -2.19%  "

# Row 11
$ws.Range("D11").Value = "3.725.90"
$ws.Range("E11").Value = "  +1.53%  "

# Row 12
$ws.Range("E12").Value = "  +18.41%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.94"
$ws.Range("E13").Value = "  -1.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.208"
$ws.Range("E14").Value = "  +1.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.91"
$ws.Range("E15").Value = "  +1.90%  "

# Row 16
$ws.Range("D16").Value = "4.424.03"
$ws.Range("E16").Value = "  +1.52%  "

# Row 17
$ws.Range("D17").Value = "97.035.93"
$ws.Range("E17").Value = "  +0.40%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.03"
$ws.Range("E18").Value = "  +1.09%  "

# Row 19
$ws.Range("D19").Value = "3.730.51"
$ws.Range("E19").Value = "  +1.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.07"
$ws.Range("E20").Value = "  +2.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.87"
$ws.Range("E21").Value = "  +0.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.505"
$ws.Range("E22").Value = "  -4.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "528.06"
$ws.Range("E23").Value = "  -0.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.48"
$ws.Range("E24").Value = "  -0.69%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000226"
$ws.Range("E25").Value = "  +10.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("E26").Value = "  -3.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "107.20"
$ws.Range("E27").Value = "  +4.61%  "

# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.934.99"
$ws.Range("E28").Value = "  +1.73%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.59"
$ws.Range("E29").Value = "  +1.09%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.190"
$ws.Range("E30").Value = "  +13.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.78"
$ws.Range("E31").Value = "  +3.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.04"
$ws.Range("E32").Value = "  +0.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.192"
$ws.Range("E34").Value = "  +2.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("E35").Value = "  -4.18%  "

# Row 36
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.31%  "

# Row 37
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.57"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "640.65"
$ws.Range("E38").Value = "  -3.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.593"
$ws.Range("E39").Value = "  -0.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.76"
$ws.Range("E40").Value = "  -0.35%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.166"
$ws.Range("E42").Value = "  +3.68%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.76"
$ws.Range("E43").Value = "  +3.15%  "

# Row 44
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.03"
$ws.Range("E44").Value = "  +1.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.58"
$ws.Range("E45").Value = "  +6.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.974"
$ws.Range("E46").Value = "  +1.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.473"
$ws.Range("E47").Value = "  +10.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0457"
$ws.Range("E48").Value = "  -0.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.38"
$ws.Range("E49").Value = "  +1.49%  "

# Row 50
$ws.Range("E50").Value = "  -0.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.66"
